$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp recorded for the existing last row (row 24) - tiny
# floating point re-computation of the same instant (2021-05-21T18:46:26Z).
$ws.Cells.Item(24, 1).Value = 44337.78224595138

# Append the newly retrieved data row (row 25) - data pulled
# Sat May 22 18:42:35 UTC 2021.
$newRow = 25
$rowValues = @(44338.77958125138, 74560, 62755, 3441, 2104, 1483, 19390, 1458, 842, 219)
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $rowValues[$col - 1]
}
